$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 25,20

$arr[0,0] = "ECs"
$arr[0,1] = "Cd38"
$arr[0,2] = "Pecam1"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 47.95321533333333
$arr[0,7] = 143.859646
$arr[0,8] = 0.6016112859309785
$arr[0,9] = 0.6016112859309785
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 405.24646
$arr[0,13] = 1215.73938
$arr[0,14] = 0.971171031955694
$arr[0,15] = 0.9711710319556939
$arr[0,16] = 19432.87075945105
$arr[0,17] = 174895.8368350595
$arr[0,18] = 0.5842674533937804
$arr[0,19] = 0.5842674533937804

$arr[1,0] = "ECs"
$arr[1,1] = "Cd38"
$arr[1,2] = "Pecam1"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 47.95321533333333
$arr[1,7] = 143.859646
$arr[1,8] = 0.6016112859309785
$arr[1,9] = 0.6016112859309785
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 5.631177666666667
$arr[1,13] = 16.893533
$arr[1,14] = 0.01349508796612936
$arr[1,15] = 0.01349508796612936
$arr[1,16] = 270.0330752299242
$arr[1,17] = 2430.297677069318
$arr[1,18] = 0.008118797225054758
$arr[1,19] = 0.008118797225054758

$arr[2,0] = "ECs"
$arr[2,1] = "Cd38"
$arr[2,2] = "Pecam1"
$arr[2,3] = "M1"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 47.95321533333333
$arr[2,7] = 143.859646
$arr[2,8] = 0.6016112859309785
$arr[2,9] = 0.6016112859309785
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 1.958728333333333
$arr[2,13] = 5.876185
$arr[2,14] = 0.004694082255041018
$arr[2,15] = 0.004694082255041017
$arr[2,16] = 93.92732154783444
$arr[2,17] = 845.3458939305099
$arr[2,18] = 0.002824012861721014
$arr[2,19] = 0.002824012861721014

$arr[3,0] = "ECs"
$arr[3,1] = "Cd38"
$arr[3,2] = "Pecam1"
$arr[3,3] = "M2"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 47.95321533333333
$arr[3,7] = 143.859646
$arr[3,8] = 0.6016112859309785
$arr[3,9] = 0.6016112859309785
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 3.263573
$arr[3,13] = 9.790719
$arr[3,14] = 0.007821135706583939
$arr[3,15] = 0.007821135706583937
$arr[3,16] = 156.4988188250526
$arr[3,17] = 1408.489369425474
$arr[3,18] = 0.004705283509878655
$arr[3,19] = 0.004705283509878654

$arr[4,0] = "ECs"
$arr[4,1] = "Cd38"
$arr[4,2] = "Pecam1"
$arr[4,3] = "sCs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 47.95321533333333
$arr[4,7] = 143.859646
$arr[4,8] = 0.6016112859309785
$arr[4,9] = 0.6016112859309785
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 1.176160333333333
$arr[4,13] = 3.528481
$arr[4,14] = 0.002818662116551706
$arr[4,15] = 0.002818662116551706
$arr[4,16] = 56.40066973085844
$arr[4,17] = 507.606027577726
$arr[4,18] = 0.001695738940543605
$arr[4,19] = 0.001695738940543605

$arr[5,0] = "FAPs"
$arr[5,1] = "Cd38"
$arr[5,2] = "Pecam1"
$arr[5,3] = "ECs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 2.315801333333333
$arr[5,7] = 6.947404000000001
$arr[5,8] = 0.0290535724960843
$arr[5,9] = 0.0290535724960843
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 405.24646
$arr[5,13] = 1215.73938
$arr[5,14] = 0.971171031955694
$arr[5,15] = 0.9711710319556939
$arr[5,16] = 938.4702923966133
$arr[5,17] = 8446.23263156952
$arr[5,18] = 0.02821598798302176
$arr[5,19] = 0.02821598798302176

$arr[6,0] = "FAPs"
$arr[6,1] = "Cd38"
$arr[6,2] = "Pecam1"
$arr[6,3] = "FAPs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 2.315801333333333
$arr[6,7] = 6.947404000000001
$arr[6,8] = 0.0290535724960843
$arr[6,9] = 0.0290535724960843
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 5.631177666666667
$arr[6,13] = 16.893533
$arr[6,14] = 0.01349508796612936
$arr[6,15] = 0.01349508796612936
$arr[6,16] = 13.04068874870356
$arr[6,17] = 117.366198738332
$arr[6,18] = 0.0003920805165649743
$arr[6,19] = 0.0003920805165649743

$arr[7,0] = "FAPs"
$arr[7,1] = "Cd38"
$arr[7,2] = "Pecam1"
$arr[7,3] = "M1"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 2.315801333333333
$arr[7,7] = 6.947404000000001
$arr[7,8] = 0.0290535724960843
$arr[7,9] = 0.0290535724960843
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 1.958728333333333
$arr[7,13] = 5.876185
$arr[7,14] = 0.004694082255041018
$arr[7,15] = 0.004694082255041017
$arr[7,16] = 4.536025685971111
$arr[7,17] = 40.82423117374
$arr[7,18] = 0.0001363798590994171
$arr[7,19] = 0.0001363798590994171

$arr[8,0] = "FAPs"
$arr[8,1] = "Cd38"
$arr[8,2] = "Pecam1"
$arr[8,3] = "M2"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 2.315801333333333
$arr[8,7] = 6.947404000000001
$arr[8,8] = 0.0290535724960843
$arr[8,9] = 0.0290535724960843
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 3.263573
$arr[8,13] = 9.790719
$arr[8,14] = 0.007821135706583939
$arr[8,15] = 0.007821135706583937
$arr[8,16] = 7.557786704830666
$arr[8,17] = 68.020080343476
$arr[8,18] = 0.00022723193325295
$arr[8,19] = 0.0002272319332529499

$arr[9,0] = "FAPs"
$arr[9,1] = "Cd38"
$arr[9,2] = "Pecam1"
$arr[9,3] = "sCs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 2.315801333333333
$arr[9,7] = 6.947404000000001
$arr[9,8] = 0.0290535724960843
$arr[9,9] = 0.0290535724960843
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 1.176160333333333
$arr[9,13] = 3.528481
$arr[9,14] = 0.002818662116551706
$arr[9,15] = 0.002818662116551706
$arr[9,16] = 2.723753668147111
$arr[9,17] = 24.513783013324
$arr[9,18] = 0.00008189220414520142
$arr[9,19] = 0.0000818922041452014

$arr[10,0] = "M1"
$arr[10,1] = "Cd38"
$arr[10,2] = "Pecam1"
$arr[10,3] = "ECs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 14.48034066666667
$arr[10,7] = 43.441022
$arr[10,8] = 0.1816674087156862
$arr[10,9] = 0.1816674087156862
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 405.24646
$arr[10,13] = 1215.73938
$arr[10,14] = 0.971171031955694
$arr[10,15] = 0.9711710319556939
$arr[10,16] = 5868.106794760708
$arr[10,17] = 52812.96115284636
$arr[10,18] = 0.1764301247951298
$arr[10,19] = 0.1764301247951298

$arr[11,0] = "M1"
$arr[11,1] = "Cd38"
$arr[11,2] = "Pecam1"
$arr[11,3] = "FAPs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 14.48034066666667
$arr[11,7] = 43.441022
$arr[11,8] = 0.1816674087156862
$arr[11,9] = 0.1816674087156862
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 5.631177666666667
$arr[11,13] = 16.893533
$arr[11,14] = 0.01349508796612936
$arr[11,15] = 0.01349508796612936
$arr[11,16] = 81.54137096785846
$arr[11,17] = 733.8723387107261
$arr[11,18] = 0.002451617661196961
$arr[11,19] = 0.002451617661196961

$arr[12,0] = "M1"
$arr[12,1] = "Cd38"
$arr[12,2] = "Pecam1"
$arr[12,3] = "M1"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 14.48034066666667
$arr[12,7] = 43.441022
$arr[12,8] = 0.1816674087156862
$arr[12,9] = 0.1816674087156862
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 1.958728333333333
$arr[12,13] = 5.876185
$arr[12,14] = 0.004694082255041018
$arr[12,15] = 0.004694082255041017
$arr[12,16] = 28.36305354011889
$arr[12,17] = 255.26748186107
$arr[12,18] = 0.0008527617595715867
$arr[12,19] = 0.0008527617595715863

$arr[13,0] = "M1"
$arr[13,1] = "Cd38"
$arr[13,2] = "Pecam1"
$arr[13,3] = "M2"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 14.48034066666667
$arr[13,7] = 43.441022
$arr[13,8] = 0.1816674087156862
$arr[13,9] = 0.1816674087156862
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 3.263573
$arr[13,13] = 9.790719
$arr[13,14] = 0.007821135706583939
$arr[13,15] = 0.007821135706583937
$arr[13,16] = 47.25764883053533
$arr[13,17] = 425.318839474818
$arr[13,18] = 0.001420845457028832
$arr[13,19] = 0.001420845457028831

$arr[14,0] = "M1"
$arr[14,1] = "Cd38"
$arr[14,2] = "Pecam1"
$arr[14,3] = "sCs"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 14.48034066666667
$arr[14,7] = 43.441022
$arr[14,8] = 0.1816674087156862
$arr[14,9] = 0.1816674087156862
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 1.176160333333333
$arr[14,13] = 3.528481
$arr[14,14] = 0.002818662116551706
$arr[14,15] = 0.002818662116551706
$arr[14,16] = 17.03120230528689
$arr[14,17] = 153.280820747582
$arr[14,18] = 0.00051205904275902
$arr[14,19] = 0.0005120590427590198

$arr[15,0] = "M2"
$arr[15,1] = "Cd38"
$arr[15,2] = "Pecam1"
$arr[15,3] = "ECs"
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 14.91142166666667
$arr[15,7] = 44.734265
$arr[15,8] = 0.1870756632601971
$arr[15,9] = 0.1870756632601971
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 405.24646
$arr[15,13] = 1215.73938
$arr[15,14] = 0.971171031955694
$arr[15,15] = 0.9711710319556939
$arr[15,16] = 6042.800843983967
$arr[15,17] = 54385.2075958557
$arr[15,18] = 0.1816824649422016
$arr[15,19] = 0.1816824649422015

$arr[16,0] = "M2"
$arr[16,1] = "Cd38"
$arr[16,2] = "Pecam1"
$arr[16,3] = "FAPs"
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 14.91142166666667
$arr[16,7] = 44.734265
$arr[16,8] = 0.1870756632601971
$arr[16,9] = 0.1870756632601971
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 5.631177666666667
$arr[16,13] = 16.893533
$arr[16,14] = 0.01349508796612936
$arr[16,15] = 0.01349508796612936
$arr[16,16] = 83.96886466758278
$arr[16,17] = 755.719782008245
$arr[16,18] = 0.002524602532018355
$arr[16,19] = 0.002524602532018355

$arr[17,0] = "M2"
$arr[17,1] = "Cd38"
$arr[17,2] = "Pecam1"
$arr[17,3] = "M1"
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 14.91142166666667
$arr[17,7] = 44.734265
$arr[17,8] = 0.1870756632601971
$arr[17,9] = 0.1870756632601971
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 1.958728333333333
$arr[17,13] = 5.876185
$arr[17,14] = 0.004694082255041018
$arr[17,15] = 0.004694082255041017
$arr[17,16] = 29.20742410878055
$arr[17,17] = 262.866816979025
$arr[17,18] = 0.0008781485512597204
$arr[17,19] = 0.0008781485512597201

$arr[18,0] = "M2"
$arr[18,1] = "Cd38"
$arr[18,2] = "Pecam1"
$arr[18,3] = "M2"
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 14.91142166666667
$arr[18,7] = 44.734265
$arr[18,8] = 0.1870756632601971
$arr[18,9] = 0.1870756632601971
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = 3.263573
$arr[18,13] = 9.790719
$arr[18,14] = 0.007821135706583939
$arr[18,15] = 0.007821135706583937
$arr[18,16] = 48.66451314294833
$arr[18,17] = 437.980618286535
$arr[18,18] = 0.001463144149757201
$arr[18,19] = 0.001463144149757201

$arr[19,0] = "M2"
$arr[19,1] = "Cd38"
$arr[19,2] = "Pecam1"
$arr[19,3] = "sCs"
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 14.91142166666667
$arr[19,7] = 44.734265
$arr[19,8] = 0.1870756632601971
$arr[19,9] = 0.1870756632601971
$arr[19,10] = 3
$arr[19,11] = 1
$arr[19,12] = 1.176160333333333
$arr[19,13] = 3.528481
$arr[19,14] = 0.002818662116551706
$arr[19,15] = 0.002818662116551706
$arr[19,16] = 17.53822267794056
$arr[19,17] = 157.844004101465
$arr[19,18] = 0.0005273030849603015
$arr[19,19] = 0.0005273030849603014

$arr[20,0] = "sCs"
$arr[20,1] = "Cd38"
$arr[20,2] = "Pecam1"
$arr[20,3] = "ECs"
$arr[20,4] = 2
$arr[20,5] = 0.6666666666666666
$arr[20,6] = 0.04719266666666667
$arr[20,7] = 0.141578
$arr[20,8] = 0.0005920695970538957
$arr[20,9] = 0.0005920695970538957
$arr[20,10] = 3
$arr[20,11] = 1
$arr[20,12] = 405.24646
$arr[20,13] = 1215.73938
$arr[20,14] = 0.971171031955694
$arr[20,15] = 0.9711710319556939
$arr[20,16] = 19.12466110462667
$arr[20,17] = 172.12194994164
$arr[20,18] = 0.0005750008415604238
$arr[20,19] = 0.0005750008415604237

$arr[21,0] = "sCs"
$arr[21,1] = "Cd38"
$arr[21,2] = "Pecam1"
$arr[21,3] = "FAPs"
$arr[21,4] = 2
$arr[21,5] = 0.6666666666666666
$arr[21,6] = 0.04719266666666667
$arr[21,7] = 0.141578
$arr[21,8] = 0.0005920695970538957
$arr[21,9] = 0.0005920695970538957
$arr[21,10] = 3
$arr[21,11] = 1
$arr[21,12] = 5.631177666666667
$arr[21,13] = 16.893533
$arr[21,14] = 0.01349508796612936
$arr[21,15] = 0.01349508796612936
$arr[21,16] = 0.2657502905637778
$arr[21,17] = 2.391752615074
$arr[21,18] = 0.000007990031294313088
$arr[21,19] = 0.000007990031294313088

$arr[22,0] = "sCs"
$arr[22,1] = "Cd38"
$arr[22,2] = "Pecam1"
$arr[22,3] = "M1"
$arr[22,4] = 2
$arr[22,5] = 0.6666666666666666
$arr[22,6] = 0.04719266666666667
$arr[22,7] = 0.141578
$arr[22,8] = 0.0005920695970538957
$arr[22,9] = 0.0005920695970538957
$arr[22,10] = 3
$arr[22,11] = 1
$arr[22,12] = 1.958728333333333
$arr[22,13] = 5.876185
$arr[22,14] = 0.004694082255041018
$arr[22,15] = 0.004694082255041017
$arr[22,16] = 0.09243761332555556
$arr[22,17] = 0.83193851993
$arr[22,18] = 0.000002779223389279977
$arr[22,19] = 0.000002779223389279977

$arr[23,0] = "sCs"
$arr[23,1] = "Cd38"
$arr[23,2] = "Pecam1"
$arr[23,3] = "M2"
$arr[23,4] = 2
$arr[23,5] = 0.6666666666666666
$arr[23,6] = 0.04719266666666667
$arr[23,7] = 0.141578
$arr[23,8] = 0.0005920695970538957
$arr[23,9] = 0.0005920695970538957
$arr[23,10] = 3
$arr[23,11] = 1
$arr[23,12] = 3.263573
$arr[23,13] = 9.790719
$arr[23,14] = 0.007821135706583939
$arr[23,15] = 0.007821135706583937
$arr[23,16] = 0.1540167127313333
$arr[23,17] = 1.386150414582
$arr[23,18] = 0.000004630656666300988
$arr[23,19] = 0.000004630656666300987

$arr[24,0] = "sCs"
$arr[24,1] = "Cd38"
$arr[24,2] = "Pecam1"
$arr[24,3] = "sCs"
$arr[24,4] = 2
$arr[24,5] = 0.6666666666666666
$arr[24,6] = 0.04719266666666667
$arr[24,7] = 0.141578
$arr[24,8] = 0.0005920695970538957
$arr[24,9] = 0.0005920695970538957
$arr[24,10] = 3
$arr[24,11] = 1
$arr[24,12] = 1.176160333333333
$arr[24,13] = 3.528481
$arr[24,14] = 0.002818662116551706
$arr[24,15] = 0.002818662116551706
$arr[24,16] = 0.05550614255755555
$arr[24,17] = 0.499555283018
$arr[24,18] = 0.000001668844143577849
$arr[24,19] = 0.000001668844143577849

$ws.Range("A2:T26").Value2 = $arr
